$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.586.61"
$ws.Range("D3").Value = "3.286.51"
$ws.Range("E3").Value = "  +4.69%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.85"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.63"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.285.38"
$ws.Range("E8").Value = "  +4.81%  "
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  +2.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.43"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +4.01%  "
$ws.Range("E12").Value = "  +2.65%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000246"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.47"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "3.834.12"
$ws.Range("E15").Value = "  +4.76%  "
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "3.288.50"
$ws.Range("E17").Value = "  +4.72%  "
$ws.Range("D18").Value = "63.686.10"
$ws.Range("E18").Value = "  +1.04%  "
$ws.Range("E19").Value = "  +2.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "477.74"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.727"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +4.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.04"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +5.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.57"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.21"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.52%  "
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.20"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.07"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.14"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.57"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +7.01%  "
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("E34").Value = "  +0.49%  "
$ws.Range("E35").Value = "  +3.75%  "
$ws.Range("E36").Value = "  +3.58%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "53.33"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.16%  "
$ws.Range("D38").Value = "0.0₃0732"
$ws.Range("E38").Value = "  +6.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0398"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +3.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "427.83"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.28%  "
$ws.Range("D41").Value = "3.072.32"
$ws.Range("E41").Value = "  +5.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.32"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("E43").Value = "  -0.22%  "
$ws.Range("E44").Value = "  +1.50%  "
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("E46").Value = "  +3.29%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.19"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.32%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "127.74"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +6.11%  "
$ws.Range("B49").Value = "USDe"
$ws.Range("C49").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.999"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "35.59"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +12.62%  "
